# Auto-generated edit script: apply numeric cell updates across 8 sheets
# per the Anima_Profits diff (commit: "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 125000780
$ws.Range("I55").Value = 200000670
$ws.Range("J55").Value = 968
$ws.Range("K55").Value = 200000670
$ws.Range("L55").Value = 968
$ws.Range("M55").Value = -200000456
$ws.Range("N55").Value = -1396
$ws.Range("H100").Value = 2082.2307
$ws.Range("I100").Value = 1915.3636
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1915.3636
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1374.3636
$ws.Range("N100").Value = -4082
$ws.Range("H116").Value = 3500
$ws.Range("I116").Value = 3285.7144
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3285.7144
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 156.2856000000002
$ws.Range("N116").Value = -11884
$ws.Range("H129").Value = 1196.5294
$ws.Range("J129").Value = 1684.0322
$ws.Range("L129").Value = 5052.096600000001
$ws.Range("N129").Value = -15052.0966
$ws.Range("H141").Value = 4448.647
$ws.Range("I141").Value = 1405
$ws.Range("K141").Value = 4215
$ws.Range("M141").Value = 965

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 959.8421
$ws.Range("I97").Value = 865.3125
$ws.Range("K97").Value = 865.3125
$ws.Range("M97").Value = -369.3125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2564.1
$ws.Range("I20").Value = 2195.5
$ws.Range("J20").Value = 2932.7
$ws.Range("K20").Value = 2195.5
$ws.Range("L20").Value = 2932.7
$ws.Range("M20").Value = -1948.5
$ws.Range("N20").Value = -3426.7
$ws.Range("H81").Value = 40519
$ws.Range("J81").Value = 40519
$ws.Range("L81").Value = 40519
$ws.Range("N81").Value = -42641
$ws.Range("H84").Value = 40519
$ws.Range("J84").Value = 40519
$ws.Range("L84").Value = 121557
$ws.Range("N84").Value = -132165
$ws.Range("H86").Value = 3281.913
$ws.Range("I86").Value = 3499.2104
$ws.Range("J86").Value = 2249.75
$ws.Range("K86").Value = 3499.2104
$ws.Range("L86").Value = 2249.75
$ws.Range("M86").Value = -2376.2104
$ws.Range("N86").Value = -4495.75
$ws.Range("H89").Value = 3281.913
$ws.Range("I89").Value = 3499.2104
$ws.Range("J89").Value = 2249.75
$ws.Range("K89").Value = 17496.052
$ws.Range("L89").Value = 11248.75
$ws.Range("M89").Value = -11880.052
$ws.Range("N89").Value = -22480.75
$ws.Range("H99").Value = 2224
$ws.Range("I99").Value = 2086
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2086
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -588
$ws.Range("N99").Value = -5496

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 58750
$ws.Range("J28").Value = 58750
$ws.Range("L28").Value = 58750
$ws.Range("N28").Value = -59240
$ws.Range("H43").Value = 109600
$ws.Range("J43").Value = 109600
$ws.Range("L43").Value = 109600
$ws.Range("N43").Value = -109968
$ws.Range("H58").Value = 1214.6316
$ws.Range("I58").Value = 905.0476
$ws.Range("K58").Value = 905.0476
$ws.Range("M58").Value = -702.0476
$ws.Range("H60").Value = 10448
$ws.Range("J60").Value = 10448
$ws.Range("L60").Value = 10448
$ws.Range("N60").Value = -11470
$ws.Range("H99").Value = 1993.6842
$ws.Range("I99").Value = 1940
$ws.Range("K99").Value = 1940
$ws.Range("M99").Value = -442
$ws.Range("H101").Value = 109600
$ws.Range("J101").Value = 109600
$ws.Range("L101").Value = 109600
$ws.Range("N101").Value = -116090
$ws.Range("H126").Value = 1993.6842
$ws.Range("I126").Value = 1940
$ws.Range("K126").Value = 5820
$ws.Range("M126").Value = -3350
$ws.Range("H129").Value = 42000
$ws.Range("J129").Value = 42000
$ws.Range("L129").Value = 42000
$ws.Range("N129").Value = -52000
$ws.Range("H132").Value = 2126.48
$ws.Range("I132").Value = 2079.1428
$ws.Range("J132").Value = 2375
$ws.Range("K132").Value = 6237.428400000001
$ws.Range("L132").Value = 7125
$ws.Range("M132").Value = -3707.428400000001
$ws.Range("N132").Value = -12185
$ws.Range("H136").Value = 1214.6316
$ws.Range("I136").Value = 905.0476
$ws.Range("K136").Value = 2715.1428
$ws.Range("M136").Value = -165.1428000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 2028.5714
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H64").Value = 1716.963
$ws.Range("I64").Value = 1024.5
$ws.Range("J64").Value = 1914.8096
$ws.Range("K64").Value = 3073.5
$ws.Range("L64").Value = 5744.4288
$ws.Range("M64").Value = -2803.5
$ws.Range("N64").Value = -6284.4288
$ws.Range("H67").Value = 1716.963
$ws.Range("I67").Value = 1024.5
$ws.Range("J67").Value = 1914.8096
$ws.Range("K67").Value = 3073.5
$ws.Range("L67").Value = 5744.4288
$ws.Range("M67").Value = -2137.5
$ws.Range("N67").Value = -7616.4288
$ws.Range("H113").Value = 550.72
$ws.Range("I113").Value = 570.86664
$ws.Range("K113").Value = 1712.59992
$ws.Range("M113").Value = 457.4000800000001
$ws.Range("H121").Value = 1127.8948
$ws.Range("I121").Value = 389.72726
$ws.Range("J121").Value = 1428.6296
$ws.Range("K121").Value = 1169.18178
$ws.Range("L121").Value = 4285.8888
$ws.Range("M121").Value = 140.8182200000001
$ws.Range("N121").Value = -6905.8888
$ws.Range("H122").Value = 3913.484
$ws.Range("I122").Value = 468.6842
$ws.Range("K122").Value = 4218.1578
$ws.Range("M122").Value = -1768.1578
$ws.Range("H131").Value = 2743.2942
$ws.Range("I131").Value = 477
$ws.Range("J131").Value = 3089
$ws.Range("K131").Value = 1431
$ws.Range("L131").Value = 9267
$ws.Range("M131").Value = 3609
$ws.Range("N131").Value = -19347

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 763.3333
$ws.Range("I97").Value = 763.3333
$ws.Range("K97").Value = 763.3333
$ws.Range("M97").Value = -267.3333
$ws.Range("H113").Value = 127587.375
$ws.Range("I113").Value = 169116.5
$ws.Range("K113").Value = 169116.5
$ws.Range("M113").Value = -166946.5
$ws.Range("H126").Value = 1802.1428
$ws.Range("J126").Value = 1728.5
$ws.Range("L126").Value = 5185.5
$ws.Range("N126").Value = -10125.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6433
$ws.Range("I7").Value = 7201.6
$ws.Range("J7").Value = 5792.5
$ws.Range("K7").Value = 7201.6
$ws.Range("L7").Value = 5792.5
$ws.Range("M7").Value = -7089.6
$ws.Range("N7").Value = -6016.5
$ws.Range("H61").Value = 3427.2222
$ws.Range("I61").Value = 2149.1667
$ws.Range("J61").Value = 5983.3335
$ws.Range("K61").Value = 2149.1667
$ws.Range("L61").Value = 5983.3335
$ws.Range("M61").Value = -1947.1667
$ws.Range("N61").Value = -6387.3335
$ws.Range("H93").Value = 18924.166
$ws.Range("I93").Value = 100000
$ws.Range("J93").Value = 2709
$ws.Range("K93").Value = 100000
$ws.Range("L93").Value = 2709
$ws.Range("M93").Value = -98752
$ws.Range("N93").Value = -5205
$ws.Range("H100").Value = 3580.4285
$ws.Range("I100").Value = 3965.75
$ws.Range("J100").Value = 3066.6667
$ws.Range("K100").Value = 3965.75
$ws.Range("L100").Value = 3066.6667
$ws.Range("M100").Value = -3424.75
$ws.Range("N100").Value = -4148.6667
$ws.Range("H113").Value = 3427.2222
$ws.Range("I113").Value = 2149.1667
$ws.Range("J113").Value = 5983.3335
$ws.Range("K113").Value = 2149.1667
$ws.Range("L113").Value = 5983.3335
$ws.Range("M113").Value = 20.83329999999978
$ws.Range("N113").Value = -10323.3335
$ws.Range("H126").Value = 6433
$ws.Range("I126").Value = 7201.6
$ws.Range("J126").Value = 5792.5
$ws.Range("K126").Value = 21604.8
$ws.Range("L126").Value = 17377.5
$ws.Range("M126").Value = -19134.8
$ws.Range("N126").Value = -22317.5
$ws.Range("H132").Value = 2888.7556
$ws.Range("I132").Value = 2953
$ws.Range("J132").Value = 2782.9412
$ws.Range("K132").Value = 8859
$ws.Range("L132").Value = 8348.8236
$ws.Range("M132").Value = -6329
$ws.Range("N132").Value = -13408.8236

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 27527.46
$ws.Range("J123").Value = 36309.5
$ws.Range("L123").Value = 36309.5
$ws.Range("N123").Value = -46109.5
$ws.Range("H132").Value = 3708.077
$ws.Range("I132").Value = 5384.2856
$ws.Range("J132").Value = 1752.5
$ws.Range("K132").Value = 16152.8568
$ws.Range("L132").Value = 5257.5
$ws.Range("M132").Value = -13622.8568
$ws.Range("N132").Value = -10317.5
$ws.Range("H136").Value = 3150.1555
$ws.Range("I136").Value = 2843.2122
$ws.Range("J136").Value = 3994.25
$ws.Range("K136").Value = 8529.6366
$ws.Range("L136").Value = 11982.75
$ws.Range("M136").Value = -5979.6366
$ws.Range("N136").Value = -17082.75
$ws.Range("H138").Value = 46159.6
$ws.Range("J138").Value = 50266
$ws.Range("L138").Value = 50266
$ws.Range("N138").Value = -60546
